$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (this also rebuilds the shared-strings table in the
#     correct order: unused old strings drop out, new ones are appended in
#     the order the cells are written) ---
$ws.Range("A2").Value = "sushanthost@test.com"
$ws.Range("C2").Value = "sushantcohost@test.com"
$ws.Range("E2").Value = "sushantguest1@test.com"
$ws.Range("L2").Value = "virtual_cabitest19"
$ws.Range("G2").Value = "michigan@na.com"

# --- Rebuild the hyperlinks. The existing ones must be removed first since
#     this object model only supports wiping the whole collection at once,
#     then every link (including the three that already existed) is
#     re-created in the original order so the relationship ids line up as
#     rId1..rId3 (existing) followed by rId4..rId6 (new). ---
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:michigan@na.com")

$ws.Hyperlinks.Add($ws.Range("K2"), "mailto:C@bi`$ush5", [Type]::Missing, [Type]::Missing, "C@bi`$ush5")
$ws.Range("K2").Value = "cabiautomation"

$ws.Hyperlinks.Add($ws.Range("I2"), "https://pushofbiz.cliotest.com/login.php")

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:sushanthost@test.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:sushantcohost@test.com")
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:sushantguest1@test.com")

# --- Re-apply the Hyperlink style to every linked cell (Hyperlinks.Add
#     leaves behind a slightly different style entry) so they all end up
#     sharing the workbook's single "Hyperlink" cell style. ---
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("G2").Style = "Hyperlink"
$ws.Range("I2").Style = "Hyperlink"
$ws.Range("K2").Style = "Hyperlink"

# --- Update the active selection shown when the sheet is opened ---
$ws.Range("G2").Select() | Out-Null
